$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 27
$ws.Range("B2").Value = "{'anOptimizer': 'adam', 'batch_size': 15, 'epochs': 23, 'hidUnit': 128, 'outActivation': 'sigmoid'}"
$ws.Range("C2").Value = 0.9632352789243063
$ws.Range("D2").Value = 0.01039861095410281
$ws.Range("E2").Value = 1

# Row 4
$ws.Range("A4").Value = 5
$ws.Range("B4").Value = "{'anOptimizer': 'rmsprop', 'batch_size': 20, 'epochs': 20, 'hidUnit': 256, 'outActivation': 'sigmoid'}"
$ws.Range("C4").Value = 0.9607843160629272
$ws.Range("D4").Value = 0.009170716395043688
$ws.Range("E4").Value = 2

# Row 5
$ws.Range("A5").Value = 21
$ws.Range("B5").Value = "{'anOptimizer': 'adam', 'batch_size': 20, 'epochs': 20, 'hidUnit': 256, 'outActivation': 'sigmoid'}"
$ws.Range("C5").Value = 0.9607843160629272
$ws.Range("D5").Value = 0.01834144872006417
$ws.Range("E5").Value = 2

# Row 6
$ws.Range("A6").Value = 23
$ws.Range("B6").Value = "{'anOptimizer': 'adam', 'batch_size': 20, 'epochs': 20, 'hidUnit': 128, 'outActivation': 'sigmoid'}"
$ws.Range("C6").Value = 0.9583333333333334
$ws.Range("D6").Value = 0.009170711085026055
$ws.Range("E6").Value = 5

# Row 7
$ws.Range("A7").Value = 19
$ws.Range("B7").Value = "{'anOptimizer': 'adam', 'batch_size': 20, 'epochs': 23, 'hidUnit': 128, 'outActivation': 'sigmoid'}"
$ws.Range("C7").Value = 0.9583333134651184
$ws.Range("D7").Value = 0.02272945092152168
$ws.Range("E7").Value = 6

# Row 8
$ws.Range("A8").Value = 29
$ws.Range("B8").Value = "{'anOptimizer': 'adam', 'batch_size': 15, 'epochs': 20, 'hidUnit': 256, 'outActivation': 'sigmoid'}"
$ws.Range("C8").Value = 0.9558823704719543
$ws.Range("D8").Value = 0.01039863905200182
$ws.Range("E8").Value = 7

# Row 9
$ws.Range("A9").Value = 17
$ws.Range("B9").Value = "{'anOptimizer': 'adam', 'batch_size': 20, 'epochs': 23, 'hidUnit': 256, 'outActivation': 'sigmoid'}"
$ws.Range("C9").Value = 0.9558823704719543
$ws.Range("D9").Value = 0.006003657055879005
$ws.Range("E9").Value = 7

# Row 10
$ws.Range("A10").Value = 25
$ws.Range("B10").Value = "{'anOptimizer': 'adam', 'batch_size': 15, 'epochs': 23, 'hidUnit': 256, 'outActivation': 'sigmoid'}"
$ws.Range("C10").Value = 0.9558823506037394
$ws.Range("D10").Value = 0.01039861095410281
$ws.Range("E10").Value = 9
